$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "89.643.37"
$ws.Range("E2").Value = "  -1.50%  "

# Row 3
$ws.Range("D3").Value = "3.081.75"
$ws.Range("E3").Value = "  -2.33%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "'235.94"
$ws.Range("E5").Value = "  +9.07%  "

# Row 6
$ws.Range("D6").Value = "'619.03"
$ws.Range("E6").Value = "  -1.16%  "

# Row 7
$ws.Range("D7").Value = "'1.06"
$ws.Range("E7").Value = "  -7.81%  "

# Row 8
$ws.Range("D8").Value = "'0.363"
$ws.Range("E8").Value = "  -1.26%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").Value = "3.078.76"
$ws.Range("E10").Value = "  -2.35%  "

# Row 11
$ws.Range("D11").Value = "'0.715"
$ws.Range("E11").Value = "  -6.21%  "

# Row 12
$ws.Range("D12").Value = "'0.199"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  +2.09%  "

# Row 14
$ws.Range("D14").Value = "'35.38"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15
$ws.Range("D15").Value = "89.337.13"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16
$ws.Range("D16").Value = "'5.37"
$ws.Range("E16").Value = "  -6.19%  "

# Row 17
$ws.Range("D17").Value = "3.641.52"
$ws.Range("E17").Value = "  -2.59%  "

# Row 18
$ws.Range("D18").Value = "3.091.42"
$ws.Range("E18").Value = "  -3.33%  "

# Row 19
$ws.Range("D19").Value = "'3.80"
$ws.Range("E19").Value = "  +0.87%  "

# Row 20
$ws.Range("D20").Value = "'0.0000213"
$ws.Range("E20").Value = "  +1.46%  "

# Row 21
$ws.Range("D21").Value = "'13.77"
$ws.Range("E21").Value = "  -5.81%  "

# Row 22
$ws.Range("D22").Value = "'433.39"
$ws.Range("E22").Value = "  -9.16%  "

# Row 23
$ws.Range("D23").Value = "'5.40"
$ws.Range("E23").Value = "  +4.29%  "

# Row 24
$ws.Range("D24").Value = "'8.78"
$ws.Range("E24").Value = "  -4.18%  "

# Row 25
$ws.Range("D25").Value = "'5.60"
$ws.Range("E25").Value = "  -5.56%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'11.75"
$ws.Range("E26").Value = "  -5.04%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'87.04"
$ws.Range("E27").Value = "  -8.21%  "

# Row 28
$ws.Range("E28").Value = "  -2.12%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.33%  "

# Row 30
$ws.Range("D30").Value = "'1.15"
$ws.Range("E30").Value = "  +14.73%  "

# Row 31
$ws.Range("D31").Value = "'9.07"
$ws.Range("E31").Value = "  -2.59%  "

# Row 32
$ws.Range("D32").Value = "'0.156"
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("D33").Value = "'0.196"
$ws.Range("E33").Value = "  -12.13%  "

# Row 34
$ws.Range("D34").Value = "'25.57"
$ws.Range("E34").Value = "  -7.01%  "

# Row 35
$ws.Range("E35").Value = "  +2.50%  "

# Row 36
$ws.Range("D36").Value = "'7.14"
$ws.Range("E36").Value = "  +2.25%  "

# Row 37
$ws.Range("D37").Value = "'3.69"
$ws.Range("E37").Value = "  +2.23%  "

# Row 38
$ws.Range("D38").Value = "'496.56"
$ws.Range("E38").Value = "  -4.58%  "

# Row 39
$ws.Range("D39").Value = "'1.88"
$ws.Range("E39").Value = "  -3.13%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").Value = "'1.26"
$ws.Range("E40").Value = "  -4.62%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0904"
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
$ws.Range("D42").Value = "'3.63"
$ws.Range("E42").Value = "  +55.92%  "

# Row 43
$ws.Range("D43").Value = "'22.08"
$ws.Range("E43").Value = "  -0.65%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").Value = "'0.398"
$ws.Range("E45").Value = "  -7.37%  "

# Row 46
$ws.Range("D46").Value = "'151.73"
$ws.Range("E46").Value = "  +1.15%  "

# Row 47
$ws.Range("E47").Value = "  -6.89%  "

# Row 48
$ws.Range("D48").Value = "'0.676"
$ws.Range("E48").Value = "  -8.23%  "

# Row 49
$ws.Range("D49").Value = "'44.49"
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.30"
$ws.Range("E50").Value = "  -5.03%  "

# Row 51
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.27%  "
